$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 530, pushing the existing 530:609 block down to 533:612
$ws.Rows("530:532").Insert()

# Row 530 - new weekly entry (Especial)
$ws.Cells.Item(530, 1).Value = 4
$ws.Cells.Item(530, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(530, 3).Value = "Los Lagos"
$ws.Cells.Item(530, 4).Value = 45127
$ws.Cells.Item(530, 5).Value = 10
$ws.Cells.Item(530, 6).Value = "Fruta"
$ws.Cells.Item(530, 7).Value = 100101
$ws.Cells.Item(530, 8).Value = "Berries"
$ws.Cells.Item(530, 9).Value = 100101007
$ws.Cells.Item(530, 10).Value = "Kiwi"
$ws.Cells.Item(530, 11).Value = "Hayward"
$ws.Cells.Item(530, 12).Value = "Especial"
$ws.Cells.Item(530, 13).Value = 200
$ws.Cells.Item(530, 14).Value = 19000
$ws.Cells.Item(530, 15).Value = 19000
$ws.Cells.Item(530, 16).Value = 19000
$ws.Cells.Item(530, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(530, 18).Value = "Regi$([char]0x00F3)n de O'Higgins"
$ws.Cells.Item(530, 19).Value = 1267
$ws.Cells.Item(530, 20).Value = 15

# Row 531 - new weekly entry (Primera)
$ws.Cells.Item(531, 1).Value = 4
$ws.Cells.Item(531, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(531, 3).Value = "Los Lagos"
$ws.Cells.Item(531, 4).Value = 45127
$ws.Cells.Item(531, 5).Value = 10
$ws.Cells.Item(531, 6).Value = "Fruta"
$ws.Cells.Item(531, 7).Value = 100101
$ws.Cells.Item(531, 8).Value = "Berries"
$ws.Cells.Item(531, 9).Value = 100101007
$ws.Cells.Item(531, 10).Value = "Kiwi"
$ws.Cells.Item(531, 11).Value = "Hayward"
$ws.Cells.Item(531, 12).Value = "Primera"
$ws.Cells.Item(531, 13).Value = 200
$ws.Cells.Item(531, 14).Value = 17000
$ws.Cells.Item(531, 15).Value = 17000
$ws.Cells.Item(531, 16).Value = 17000
$ws.Cells.Item(531, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(531, 18).Value = "Regi$([char]0x00F3)n de O'Higgins"
$ws.Cells.Item(531, 19).Value = 1133
$ws.Cells.Item(531, 20).Value = 15

# Row 532 - new weekly entry (Segunda)
$ws.Cells.Item(532, 1).Value = 4
$ws.Cells.Item(532, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(532, 3).Value = "Los Lagos"
$ws.Cells.Item(532, 4).Value = 45127
$ws.Cells.Item(532, 5).Value = 10
$ws.Cells.Item(532, 6).Value = "Fruta"
$ws.Cells.Item(532, 7).Value = 100101
$ws.Cells.Item(532, 8).Value = "Berries"
$ws.Cells.Item(532, 9).Value = 100101007
$ws.Cells.Item(532, 10).Value = "Kiwi"
$ws.Cells.Item(532, 11).Value = "Hayward"
$ws.Cells.Item(532, 12).Value = "Segunda"
$ws.Cells.Item(532, 13).Value = 200
$ws.Cells.Item(532, 14).Value = 13000
$ws.Cells.Item(532, 15).Value = 13000
$ws.Cells.Item(532, 16).Value = 13000
$ws.Cells.Item(532, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(532, 18).Value = "Regi$([char]0x00F3)n de O'Higgins"
$ws.Cells.Item(532, 19).Value = 867
$ws.Cells.Item(532, 20).Value = 15

# Make sure the date column keeps the date-time number format used elsewhere in column D
$ws.Range("D530:D532").NumberFormat = $ws.Range("D533").NumberFormat
